$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Table-of-contents block rework (paragraphs 23-27 before edit)
#    - the empty paragraph gets a "single" underline on its mark and
#      the "_GoBack" bookmark moves there
#    - each TOC line gets a dotted-leader + page-number run appended
# ---------------------------------------------------------------

$tocEmpty = $d.Paragraphs.Item(23)

# Give the (still empty) paragraph mark a "single" underline. A
# paragraph-mark-only range can't carry character formatting on its
# own, so type a placeholder character, format the whole (now
# non-empty) paragraph range, then delete the placeholder again - the
# mark keeps the formatting, exactly as it would editing live in Word.
$markStart = $tocEmpty.Range.Start
$placeholder = $d.Range($markStart, $markStart)
$placeholder.InsertAfter("X")
$tocEmpty2 = $d.Paragraphs.Item(23)
$tocEmpty2.Range.Font.Underline = 1
$charRange = $d.Range($markStart, $markStart + 1)
$charRange.Delete()

# Move the "_GoBack" bookmark onto this paragraph's mark (re-adding a
# bookmark with the same name moves it and drops the old position).
$tocEmpty3 = $d.Paragraphs.Item(23)
$d.Bookmarks.Add("_GoBack", $tocEmpty3.Range)

# Helper: append a run of text right before a paragraph's end mark,
# as a distinct run (so it matches a freshly-typed trailing run
# rather than being silently re-merged into the preceding one).
function Append-SplitRun($paraIndex, $text) {
    $para = $d.Paragraphs.Item($paraIndex)
    $full = $para.Range
    $body = $d.Range($full.Start, $full.End - 1)
    $before = $body.End
    $body.InsertAfter($text)
    $newRange = $d.Range($before, $body.End)
    $newRange.Font.Bold = 1
    $newRange.Font.Bold = 0
}

# Plain append - used where a bookmark (or other element) between the
# previous run and the new text already forces a fresh <w:r>.
function Append-Plain($paraIndex, $text) {
    $para = $d.Paragraphs.Item($paraIndex)
    $full = $para.Range
    $body = $d.Range($full.Start, $full.End - 1)
    $body.InsertAfter($text)
}

Append-SplitRun 24 " ………………………………………………………………. 3"
Append-SplitRun 25 " ……………………………………………………………………………………. 3"
Append-SplitRun 26 " ……………………………………………………………. 4"
Append-Plain    27 " ……………………………………………………… 4"
